$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new motor entry as row 11
$ws.Range("A11").Value = "SunnySky X Series V3 X2305 KV1450"
$ws.Range("B11").Value = 1450
$ws.Range("C11").Value = 0.12
$ws.Range("D11").Value = 20.5
$ws.Range("E11").Value = 28
$ws.Range("F11").Value = 20
$ws.Range("G11").Value = 12.6
$ws.Range("H11").Value = 25

# Update selection to reflect the next empty row (A12), matching source workbook state
$ws.Range("A12").Select()
